$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 3820.25
$ws.Range("I100").Value = 4278.1113
$ws.Range("J100").Value = 3231.5715
$ws.Range("K100").Value = 4278.1113
$ws.Range("L100").Value = 3231.5715
$ws.Range("M100").Value = -3737.1113
$ws.Range("N100").Value = -4313.5715
$ws.Range("H113").Value = 3076.875
$ws.Range("I113").Value = 2024
$ws.Range("J113").Value = 3427.8333
$ws.Range("K113").Value = 2024
$ws.Range("L113").Value = 3427.8333
$ws.Range("M113").Value = 1230
$ws.Range("N113").Value = -9935.8333
$ws.Range("H132").Value = 1303.8784
$ws.Range("I132").Value = 1351.2715
$ws.Range("K132").Value = 4053.8145
$ws.Range("M132").Value = -1523.8145
$ws.Range("H138").Value = 2016.6
$ws.Range("J138").Value = 4685.2144
$ws.Range("L138").Value = 14055.6432
$ws.Range("N138").Value = -24335.6432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 100000
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H45").Value = 9531.23
$ws.Range("I45").Value = 10264.363
$ws.Range("K45").Value = 10264.363
$ws.Range("M45").Value = -9887.362999999999
$ws.Range("H61").Value = 13733.56
$ws.Range("I61").Value = 2036.6885
$ws.Range("K61").Value = 2036.6885
$ws.Range("M61").Value = -1824.6885
$ws.Range("H74").Value = 66785.53999999999
$ws.Range("I74").Value = 42955.047
$ws.Range("K74").Value = 42955.047
$ws.Range("M74").Value = -42081.047
$ws.Range("H77").Value = 66785.53999999999
$ws.Range("I77").Value = 42955.047
$ws.Range("K77").Value = 214775.235
$ws.Range("M77").Value = -210407.235
$ws.Range("H110").Value = 183841
$ws.Range("I110").Value = 301985.16
$ws.Range("K110").Value = 301985.16
$ws.Range("M110").Value = -299940.16
$ws.Range("H119").Value = 80597.8
$ws.Range("J119").Value = 80597.8
$ws.Range("L119").Value = 80597.8
$ws.Range("N119").Value = -90273.8
$ws.Range("H132").Value = 2396.75
$ws.Range("I132").Value = 2104.1833
$ws.Range("J132").Value = 3274.45
$ws.Range("K132").Value = 6312.5499
$ws.Range("L132").Value = 9823.349999999999
$ws.Range("M132").Value = -3782.5499
$ws.Range("N132").Value = -14883.35
$ws.Range("H136").Value = 13733.56
$ws.Range("I136").Value = 2036.6885
$ws.Range("K136").Value = 6110.0655
$ws.Range("M136").Value = -3560.0655

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 2467.111
$ws.Range("I5").Value = 386.7143
$ws.Range("K5").Value = 386.7143
$ws.Range("M5").Value = -273.7143
$ws.Range("H22").Value = 230
$ws.Range("I22").Value = 230
$ws.Range("K22").Value = 230
$ws.Range("M22").Value = -57
$ws.Range("H99").Value = 1347.3125
$ws.Range("I99").Value = 1327.5385
$ws.Range("J99").Value = 1433
$ws.Range("K99").Value = 1327.5385
$ws.Range("L99").Value = 1433
$ws.Range("M99").Value = 170.4614999999999
$ws.Range("N99").Value = -4429
$ws.Range("H107").Value = 6796.196
$ws.Range("I107").Value = 6148.0527
$ws.Range("K107").Value = 6148.0527
$ws.Range("M107").Value = -4228.0527

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2409.0967
$ws.Range("I16").Value = 1954.7826
$ws.Range("J16").Value = 3715.25
$ws.Range("K16").Value = 1954.7826
$ws.Range("L16").Value = 3715.25
$ws.Range("M16").Value = -1667.7826
$ws.Range("N16").Value = -4289.25
$ws.Range("H22").Value = 576.76666
$ws.Range("I22").Value = 427.2857
$ws.Range("K22").Value = 427.2857
$ws.Range("M22").Value = -77.28570000000002
$ws.Range("H58").Value = 3905.4211
$ws.Range("I58").Value = 1904.1702
$ws.Range("J58").Value = 13311.3
$ws.Range("K58").Value = 1904.1702
$ws.Range("L58").Value = 13311.3
$ws.Range("M58").Value = -1701.1702
$ws.Range("N58").Value = -13717.3
$ws.Range("H88").Value = 23624.375
$ws.Range("J88").Value = 22383.428
$ws.Range("L88").Value = 22383.428
$ws.Range("N88").Value = -23195.428
$ws.Range("H91").Value = 23624.375
$ws.Range("J91").Value = 22383.428
$ws.Range("L91").Value = 22383.428
$ws.Range("N91").Value = -25191.428
$ws.Range("H113").Value = 2409.0967
$ws.Range("I113").Value = 1954.7826
$ws.Range("J113").Value = 3715.25
$ws.Range("K113").Value = 1954.7826
$ws.Range("L113").Value = 3715.25
$ws.Range("M113").Value = 215.2174
$ws.Range("N113").Value = -8055.25
$ws.Range("H134").Value = 2173.75
$ws.Range("I134").Value = 2102.2712
$ws.Range("J134").Value = 2642.3333
$ws.Range("K134").Value = 6306.8136
$ws.Range("L134").Value = 7926.999899999999
$ws.Range("M134").Value = -3771.8136
$ws.Range("N134").Value = -12996.9999
$ws.Range("H136").Value = 3905.4211
$ws.Range("I136").Value = 1904.1702
$ws.Range("J136").Value = 13311.3
$ws.Range("K136").Value = 5712.5106
$ws.Range("L136").Value = 39933.89999999999
$ws.Range("M136").Value = -3162.5106
$ws.Range("N136").Value = -45033.89999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1048.1538
$ws.Range("I5").Value = 1339.5
$ws.Range("K5").Value = 4018.5
$ws.Range("M5").Value = -3906.5
$ws.Range("H135").Value = 1048.1538
$ws.Range("I135").Value = 1339.5
$ws.Range("K135").Value = 12055.5
$ws.Range("M135").Value = -9520.5
$ws.Range("H140").Value = 1846.2069
$ws.Range("I140").Value = 1503.88
$ws.Range("K140").Value = 4511.64
$ws.Range("M140").Value = 668.3599999999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 19083.652
$ws.Range("I102").Value = 21734.615
$ws.Range("K102").Value = 21734.615
$ws.Range("M102").Value = -20112.615

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5859.826
$ws.Range("I22").Value = 7177.2856
$ws.Range("J22").Value = 5283.4375
$ws.Range("K22").Value = 7177.2856
$ws.Range("L22").Value = 5283.4375
$ws.Range("M22").Value = -6882.2856
$ws.Range("N22").Value = -5873.4375
$ws.Range("H27").Value = 5859.826
$ws.Range("I27").Value = 7177.2856
$ws.Range("J27").Value = 5283.4375
$ws.Range("K27").Value = 7177.2856
$ws.Range("L27").Value = 5283.4375
$ws.Range("M27").Value = -7070.2856
$ws.Range("N27").Value = -5497.4375
$ws.Range("H40").Value = 28174.5
$ws.Range("I40").Value = 34139.07
$ws.Range("K40").Value = 34139.07
$ws.Range("M40").Value = -34003.07
$ws.Range("H55").Value = 1414.9445
$ws.Range("J55").Value = 3270.2
$ws.Range("L55").Value = 3270.2
$ws.Range("N55").Value = -3616.2
$ws.Range("H100").Value = 2371.3333
$ws.Range("I100").Value = 2045.6
$ws.Range("K100").Value = 2045.6
$ws.Range("M100").Value = -1504.6
$ws.Range("H136").Value = 25045.926
$ws.Range("J136").Value = 94464.5
$ws.Range("L136").Value = 283393.5
$ws.Range("N136").Value = -288493.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1961.4615
$ws.Range("I81").Value = 1450
$ws.Range("K81").Value = 2900
$ws.Range("M81").Value = -1839
$ws.Range("H84").Value = 1961.4615
$ws.Range("I84").Value = 1450
$ws.Range("K84").Value = 14500
$ws.Range("M84").Value = -9196
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("M119").ClearContents()
$ws.Range("N119").ClearContents()
$ws.Range("H132").Value = 2620.2632
$ws.Range("I132").Value = 2761.963
$ws.Range("K132").Value = 8285.889000000001
$ws.Range("M132").Value = -5755.889000000001
